# Append order row 12 to the "Orders" sheet.
# Values are digit-looking text ("12"), so a leading apostrophe forces
# Excel to store them as text (matching the workbook's existing
# numberStoredAsText convention) instead of coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$ws.Range("A12").Value = "'12"
$ws.Range("C12").Value = "447_黄金球_craspedia_undefined_1bunch"
$ws.Range("F12").Value = "'12"

# Update the Summary sheet's tracking/reference number in G2.
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("G2").Value = "'010135331020662812"
